$d = $word.ActiveDocument

# Locate the last paragraph of the document body (the Conclusion text),
# which ends with "...ensure high-quality collaboration."
$lastPara = $d.Paragraphs.Last

# Insert a brand new paragraph right after it.
$newRange = $lastPara.Range.InsertParagraphAfter()

# The new paragraph's range is the just-inserted paragraph mark; grab the
# paragraph object itself so we can fill in its text run.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "(Added by Eshal) " + [char]0x2014 + " Using GitHub taught me how developers manage code and collaborate effectively on real projects."

# Match formatting used by the document's body text (Times New Roman).
$newPara.Range.Font.Name = "Times New Roman"
